$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Correct the "FilesTab" Cypher query text in B4 ---
# Remove the `File Type` and `Breed` columns from the RETURN clause
# (the two coalesce(...) lines that used to project f.file_type and demo.breed).
$newFilesQuery = @"
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
WHERE demo.breed IN ['Cavalier King Charles Spaniel']
OPTIONAL MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
WITH DISTINCT f, parent, c, demo, diag, s
RETURN  coalesce(f.file_name, '') AS ``File Name``,
        coalesce(labels(parent)[0], '') AS ``Association``,
        coalesce(f.file_description, '') AS ``Description``,
        coalesce(f.file_format, '') AS ``Format``,
        coalesce(f.file_size, '') AS ``Size``,
        coalesce(c.case_id, '') AS ``Case ID``,
        coalesce(diag.disease_term,'') AS Diagnosis , 
        coalesce(s.clinical_study_designation,'') AS ``Study Code``
"@

$ws.Range("B4").Value = $newFilesQuery

# --- Row 4 grew shorter now that two lines of the query were removed ---
$ws.Range("B4").EntireRow.RowHeight = 217.5

# --- View state: window scrolled back to the top-left, zoomed out to 70%, ---
# --- and the active cell moved from E4 to B4 ---
$win = $excel.ActiveWindow
$win.Zoom = 70
$ws.Range("B4").Select()
